$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the neighboring header cell (G1) to the new H1 header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
